# Apply the edits described by the diff:
#  1. Update the "Date" metadata value on the "Metadata" sheet.
#  2. Update Min/Max (and Base Min/Base Max) for the
#     "ActorPS.XCN9.composant1" row on the "Elements" sheet from 1 to 0.

$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-05-05T11:54:16+00:00"

# --- 2. Elements sheet: row 7 Min/Max and Base Min/Base Max ---
# These columns store numeric-looking values ("0"/"1") as *text* (shared
# strings), so a leading apostrophe is used to force text entry and keep
# them from being reinterpreted as numbers.
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("F7").Value = "'0"
$elements.Range("G7").Value = "'0"
$elements.Range("AG7").Value = "'0"
$elements.Range("AH7").Value = "'0"
